$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.269.82'
$ws.Range("E2").Value = '  +1.75%  '
$ws.Range("D3").Value = '2.348.55'
$ws.Range("E3").Value = '  -1.54%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '''545.36'
$ws.Range("E5").Value = '  +1.20%  '
$ws.Range("D6").Value = '''136.67'
$ws.Range("E6").Value = '  -2.25%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").Value = '''0.526'
$ws.Range("E8").Value = '  -8.30%  '
$ws.Range("D9").Value = '2.347.51'
$ws.Range("E9").Value = '  -1.37%  '
$ws.Range("E10").Value = '  +0.67%  '
$ws.Range("D11").Value = '''0.157'
$ws.Range("E11").Value = '  +1.76%  '
$ws.Range("D12").Value = '''5.33'
$ws.Range("E12").Value = '  +0.21%  '
$ws.Range("D13").Value = '''0.341'
$ws.Range("E13").Value = '  +0.48%  '
$ws.Range("E14").Value = '  -2.31%  '
$ws.Range("D15").Value = '2.773.89'
$ws.Range("E15").Value = '  -1.51%  '
$ws.Range("D16").Value = '60.763.00'
$ws.Range("E16").Value = '  +1.10%  '
$ws.Range("E17").Value = '  -2.47%  '
$ws.Range("D18").Value = '2.348.19'
$ws.Range("E18").Value = '  -1.78%  '
$ws.Range("E19").Value = '  +0.46%  '
$ws.Range("D20").Value = '''318.78'
$ws.Range("E20").Value = '  +1.70%  '
$ws.Range("E21").Value = '  +1.52%  '
$ws.Range("D22").Value = '''6.54'
$ws.Range("E22").Value = '  -2.14%  '
$ws.Range("D23").Value = '''1.00'
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").Value = '''1.75'
$ws.Range("E24").Value = '  -2.43%  '
$ws.Range("D25").Value = '''63.22'
$ws.Range("E25").Value = '  +1.03%  '
$ws.Range("D26").Value = '''8.28'
$ws.Range("E26").Value = '  +8.36%  '
$ws.Range("D27").Value = '''7.95'
$ws.Range("E27").Value = '  -0.32%  '
$ws.Range("D28").Value = '''498.81'
$ws.Range("E28").Value = '  -1.39%  '
$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").Value = '''0.148'
$ws.Range("E29").Value = '  +3.67%  '
$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D30").Value = '''1.38'
$ws.Range("E30").Value = '  -2.43%  '
$ws.Range("B31").Value = 'PEPE'
$ws.Range("C31").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D31").Value = '0.0₃0868'
$ws.Range("E31").Value = '  -4.55%  '
$ws.Range("E32").Value = '  -1.73%  '
$ws.Range("E33").Value = '  -3.22%  '
$ws.Range("E34").Value = '  -0.06%  '
$ws.Range("D35").Value = '''4.60'
$ws.Range("E35").Value = '  -0.32%  '
$ws.Range("D36").Value = '''0.377'
$ws.Range("E36").Value = '  +1.58%  '
$ws.Range("D37").Value = '''18.50'
$ws.Range("E37").Value = '  +3.08%  '
$ws.Range("E38").Value = '  -2.46%  '
$ws.Range("E39").Value = '  +6.82%  '
$ws.Range("D40").Value = '''141.77'
$ws.Range("E40").Value = '  +3.59%  '
$ws.Range("D41").Value = '''1.00'
$ws.Range("E41").Value = '  +0.21%  '
$ws.Range("D42").Value = '''40.40'
$ws.Range("E42").Value = '  +0.34%  '
$ws.Range("D43").Value = '''142.33'
$ws.Range("E43").Value = '  +1.99%  '
$ws.Range("D44").Value = '''3.55'
$ws.Range("E44").Value = '  +1.21%  '
$ws.Range("E45").Value = '  -4.70%  '
$ws.Range("E46").Value = '  +0.31%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").Value = '''19.07'
$ws.Range("E47").Value = '  -5.38%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").Value = '''0.559'
$ws.Range("E48").Value = '  -2.64%  '
$ws.Range("D49").Value = '''0.0903'
$ws.Range("E49").Value = '  -2.01%  '
$ws.Range("D50").Value = '''0.0221'
$ws.Range("E50").Value = '  -0.72%  '
$ws.Range("E51").Value = '  -2.67%  '
